# Test Type Data type correction
# Insert a new validation-error row (for the "amount" column,
# CUSTOM_LOGIC_VIOLATION) right after the existing row for Row_Number=42,
# shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4; existing rows 4-31 shift down to 5-32.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new validation error entry.
$ws.Cells.Item(4, 1).Value = 61
$ws.Cells.Item(4, 2).Value = 2
$ws.Cells.Item(4, 3).Value = "amount"
$ws.Cells.Item(4, 4).Value = 150
$ws.Cells.Item(4, 5).Value = "CUSTOM_LOGIC_VIOLATION"
$ws.Cells.Item(4, 6).Value = "DSL Rule failure"
